$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# folder5 / folder6 (and their text children) actually live under
# folder\folder2\folder4, not folder\folder2\folder3 -- fix the path column
# for those rows (11-14) so testing of the new 'folder' fixture lines up.
$ws.Range("E11").Value = "folder\folder2\folder4"
$ws.Range("E12").Value = "folder\folder2\folder4"
$ws.Range("E13").Value = "folder\folder2\folder4"
$ws.Range("E14").Value = "folder\folder2\folder4"

# Add a new "level" header in A1 (was blank) and rename the "folders" header
# in C1 to "file_name" -- this documents the make_setup.py output columns.
$ws.Range("A1").Value = "level"
$ws.Range("C1").Value = "file_name"

# Move the active selection to E14.
$ws.Range("E14").Select()
